$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: 'Datos actualizados a 27 de Marzo de 2020 a las 16:14' -> 'Datos actualizados a 27 de Marzo de 2020 a las 16:44'
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Marzo de 2020 a las 16:44"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 86548
$ws.Cells.Item(4, 3).Value = 1113
$ws.Cells.Item(4, 4).Value = 1889
$ws.Cells.Item(4, 5).Value = 83338
$ws.Cells.Item(4, 6).Value = 2139
$ws.Cells.Item(4, 7).Value = 26
$ws.Cells.Item(4, 8).Value = 1321

# Row 7: 'España' -> 'España'
$ws.Cells.Item(7, 1).Value = "España"
$ws.Cells.Item(7, 2).Value = 64059
$ws.Cells.Item(7, 3).Value = 6273
$ws.Cells.Item(7, 4).Value = 9357
$ws.Cells.Item(7, 5).Value = 49768
$ws.Cells.Item(7, 6).Value = 4165
$ws.Cells.Item(7, 7).Value = 569
$ws.Cells.Item(7, 8).Value = 4934

# Row 11: 'Reino Unido' -> 'Reino Unido'
$ws.Cells.Item(11, 1).Value = "Reino Unido"
$ws.Cells.Item(11, 2).Value = 14543
$ws.Cells.Item(11, 3).Value = 2885
$ws.Cells.Item(11, 4).Value = 135
$ws.Cells.Item(11, 5).Value = 13649
$ws.Cells.Item(11, 6).Value = 163
$ws.Cells.Item(11, 7).Value = 181
$ws.Cells.Item(11, 8).Value = 759

# Row 15: 'Austria' -> 'Austria'
$ws.Cells.Item(15, 1).Value = "Austria"
$ws.Cells.Item(15, 2).Value = 7500
$ws.Cells.Item(15, 3).Value = 591
$ws.Cells.Item(15, 4).Value = 225
$ws.Cells.Item(15, 5).Value = 7217
$ws.Cells.Item(15, 6).Value = 128
$ws.Cells.Item(15, 7).Value = 9
$ws.Cells.Item(15, 8).Value = 58

# Row 30: 'Luxemburgo' -> 'Luxemburgo'
$ws.Cells.Item(30, 1).Value = "Luxemburgo"
$ws.Cells.Item(30, 2).Value = 1605
$ws.Cells.Item(30, 3).Value = 152
$ws.Cells.Item(30, 4).Value = 6
$ws.Cells.Item(30, 5).Value = 1584
$ws.Cells.Item(30, 6).Value = 3
$ws.Cells.Item(30, 7).Value = 6
$ws.Cells.Item(30, 8).Value = 15

# Row 44: 'India' -> 'India'
$ws.Cells.Item(44, 1).Value = "India"
$ws.Cells.Item(44, 2).Value = 874
$ws.Cells.Item(44, 3).Value = 147
$ws.Cells.Item(44, 4).Value = 73
$ws.Cells.Item(44, 5).Value = 781
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 20

# Row 70: 'Letonia' -> 'Bulgaria'
$ws.Cells.Item(70, 1).Value = "Bulgaria"
$ws.Cells.Item(70, 2).Value = 293
$ws.Cells.Item(70, 3).Value = 29
$ws.Cells.Item(70, 4).Value = 9
$ws.Cells.Item(70, 5).Value = 281
$ws.Cells.Item(70, 6).Value = 8
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 3

# Row 71: 'Bulgaria' -> 'Letonia'
$ws.Cells.Item(71, 1).Value = "Letonia"
$ws.Cells.Item(71, 2).Value = 280
$ws.Cells.Item(71, 3).Value = 36
$ws.Cells.Item(71, 4).Value = 1
$ws.Cells.Item(71, 5).Value = 279
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0

# Row 85: 'Albania' -> 'Moldavia'
$ws.Cells.Item(85, 1).Value = "Moldavia"
$ws.Cells.Item(85, 2).Value = 199
$ws.Cells.Item(85, 3).Value = 22
$ws.Cells.Item(85, 4).Value = 2
$ws.Cells.Item(85, 5).Value = 195
$ws.Cells.Item(85, 6).Value = 33
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(85, 8).Value = 2

# Row 86: 'Moldavia' -> 'Albania'
$ws.Cells.Item(86, 1).Value = "Albania"
$ws.Cells.Item(86, 2).Value = 186
$ws.Cells.Item(86, 3).Value = 12
$ws.Cells.Item(86, 4).Value = 31
$ws.Cells.Item(86, 5).Value = 147
$ws.Cells.Item(86, 6).Value = 3
$ws.Cells.Item(86, 7).Value = 2
$ws.Cells.Item(86, 8).Value = 8

# Row 87: 'Azerbaiyan' -> 'Burkina Faso'
$ws.Cells.Item(87, 1).Value = "Burkina Faso"
$ws.Cells.Item(87, 2).Value = 180
$ws.Cells.Item(87, 3).Value = 28
$ws.Cells.Item(87, 4).Value = 12
$ws.Cells.Item(87, 5).Value = 159
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 9

# Row 88: 'Vietnam' -> 'Azerbaiyan'
$ws.Cells.Item(88, 1).Value = "Azerbaiyan"
$ws.Cells.Item(88, 2).Value = 165
$ws.Cells.Item(88, 3).Value = 43
$ws.Cells.Item(88, 4).Value = 15
$ws.Cells.Item(88, 5).Value = 147
$ws.Cells.Item(88, 6).Value = 6
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 3

# Row 89: 'Burkina Faso' -> 'Vietnam'
$ws.Cells.Item(89, 1).Value = "Vietnam"
$ws.Cells.Item(89, 2).Value = 163
$ws.Cells.Item(89, 3).Value = 10
$ws.Cells.Item(89, 4).Value = 20
$ws.Cells.Item(89, 5).Value = 143
$ws.Cells.Item(89, 6).Value = 3
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 0

# Row 94: 'Ghana' -> 'Kazajistan'
$ws.Cells.Item(94, 1).Value = "Kazajistan"
$ws.Cells.Item(94, 2).Value = 137
$ws.Cells.Item(94, 3).Value = 24
$ws.Cells.Item(94, 4).Value = 3
$ws.Cells.Item(94, 5).Value = 133
$ws.Cells.Item(94, 6).Value = 0
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = 1

# Row 95: 'Kazajistan' -> 'Ghana'
$ws.Cells.Item(95, 1).Value = "Ghana"
$ws.Cells.Item(95, 2).Value = 136
$ws.Cells.Item(95, 3).Value = 4
$ws.Cells.Item(95, 4).Value = 1
$ws.Cells.Item(95, 5).Value = 131
$ws.Cells.Item(95, 6).Value = 1
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 4

# Row 108: 'Uzbekistan' -> 'Uzbekistan'
$ws.Cells.Item(108, 1).Value = "Uzbekistan"
$ws.Cells.Item(108, 2).Value = 88
$ws.Cells.Item(108, 3).Value = 13
$ws.Cells.Item(108, 4).Value = 5
$ws.Cells.Item(108, 5).Value = 82
$ws.Cells.Item(108, 6).Value = 8
$ws.Cells.Item(108, 7).Value = 1
$ws.Cells.Item(108, 8).Value = 1

# Row 111: 'Guadalupe' -> 'Cuba'
$ws.Cells.Item(111, 1).Value = "Cuba"
$ws.Cells.Item(111, 2).Value = 80
$ws.Cells.Item(111, 3).Value = 13
$ws.Cells.Item(111, 4).Value = 4
$ws.Cells.Item(111, 5).Value = 74
$ws.Cells.Item(111, 6).Value = 2
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 2

# Row 112: 'Montenegro' -> 'Guadalupe'
$ws.Cells.Item(112, 1).Value = "Guadalupe"
$ws.Cells.Item(112, 2).Value = 73
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 0
$ws.Cells.Item(112, 5).Value = 72
$ws.Cells.Item(112, 6).Value = 4
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 1

# Row 113: 'Honduras' -> 'Montenegro'
$ws.Cells.Item(113, 1).Value = "Montenegro"
$ws.Cells.Item(113, 2).Value = 70
$ws.Cells.Item(113, 3).Value = 1
$ws.Cells.Item(113, 4).Value = 0
$ws.Cells.Item(113, 5).Value = 69
$ws.Cells.Item(113, 6).Value = 1
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 1

# Row 114: 'Cuba' -> 'Honduras'
$ws.Cells.Item(114, 1).Value = "Honduras"
$ws.Cells.Item(114, 2).Value = 68
$ws.Cells.Item(114, 3).Value = 1
$ws.Cells.Item(114, 4).Value = 0
$ws.Cells.Item(114, 5).Value = 67
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 1

# Row 120: 'Paraguay' -> 'Gibraltar'
$ws.Cells.Item(120, 1).Value = "Gibraltar"
$ws.Cells.Item(120, 2).Value = 55
$ws.Cells.Item(120, 3).Value = 20
$ws.Cells.Item(120, 4).Value = 14
$ws.Cells.Item(120, 5).Value = 41
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 0

# Row 121: 'Consejo Danes para los Refugiados' -> 'Paraguay'
$ws.Cells.Item(121, 1).Value = "Paraguay"
$ws.Cells.Item(121, 2).Value = 52
$ws.Cells.Item(121, 3).Value = 11
$ws.Cells.Item(121, 4).Value = 1
$ws.Cells.Item(121, 5).Value = 48
$ws.Cells.Item(121, 6).Value = 1
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 3

# Row 122: 'Mayotte' -> 'Consejo Danes para los Refugiados'
$ws.Cells.Item(122, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(122, 2).Value = 51
$ws.Cells.Item(122, 3).Value = 0
$ws.Cells.Item(122, 4).Value = 2
$ws.Cells.Item(122, 5).Value = 46
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 3

# Row 123: 'Ruanda' -> 'Mayotte'
$ws.Cells.Item(123, 1).Value = "Mayotte"
$ws.Cells.Item(123, 2).Value = 50
$ws.Cells.Item(123, 3).Value = 14
$ws.Cells.Item(123, 4).Value = 0
$ws.Cells.Item(123, 5).Value = 50
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 0

# Row 124: 'Banglades' -> 'Ruanda'
$ws.Cells.Item(124, 1).Value = "Ruanda"
$ws.Cells.Item(124, 2).Value = 50
$ws.Cells.Item(124, 3).Value = 0
$ws.Cells.Item(124, 4).Value = 0
$ws.Cells.Item(124, 5).Value = 50
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 0

# Row 125: 'Puerto Rico' -> 'Banglades'
$ws.Cells.Item(125, 1).Value = "Banglades"
$ws.Cells.Item(125, 2).Value = 48
$ws.Cells.Item(125, 3).Value = 4
$ws.Cells.Item(125, 4).Value = 11
$ws.Cells.Item(125, 5).Value = 32
$ws.Cells.Item(125, 6).Value = 1
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 5

# Row 126: 'Gibraltar' -> 'Puerto Rico'
$ws.Cells.Item(126, 1).Value = "Puerto Rico"
$ws.Cells.Item(126, 2).Value = 39
$ws.Cells.Item(126, 3).Value = 0
$ws.Cells.Item(126, 4).Value = 1
$ws.Cells.Item(126, 5).Value = 36
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 2

# Row 138: 'Barbados' -> 'Madagascar'
$ws.Cells.Item(138, 1).Value = "Madagascar"
$ws.Cells.Item(138, 2).Value = 24
$ws.Cells.Item(138, 3).Value = 1
$ws.Cells.Item(138, 4).Value = 0
$ws.Cells.Item(138, 5).Value = 24
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 0

# Row 139: 'Madagascar' -> 'Barbados'
$ws.Cells.Item(139, 1).Value = "Barbados"
$ws.Cells.Item(139, 2).Value = 24
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(139, 4).Value = 0
$ws.Cells.Item(139, 5).Value = 24
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 0

# Row 149: 'Guinea Ecuatorial' -> 'Republica de Yibuti'
$ws.Cells.Item(149, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(149, 2).Value = 12
$ws.Cells.Item(149, 3).Value = 1
$ws.Cells.Item(149, 4).Value = 0
$ws.Cells.Item(149, 5).Value = 12
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 0

# Row 150: 'Republica de Yibuti' -> 'Guinea Ecuatorial'
$ws.Cells.Item(150, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(150, 2).Value = 12
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 0
$ws.Cells.Item(150, 5).Value = 12
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 0

# Row 151: 'San Martin (Parte Francesa)' -> 'Mongolia'
$ws.Cells.Item(151, 1).Value = "Mongolia"
$ws.Cells.Item(151, 2).Value = 11
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 0
$ws.Cells.Item(151, 5).Value = 11
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 0

# Row 152: 'Mongolia' -> 'San Martin (Parte Francesa)'
$ws.Cells.Item(152, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(152, 2).Value = 11
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = 0
$ws.Cells.Item(152, 5).Value = 11
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 0

# Row 153: 'Dominica' -> 'Mali'
$ws.Cells.Item(153, 1).Value = "Mali"
$ws.Cells.Item(153, 2).Value = 11
$ws.Cells.Item(153, 3).Value = 7
$ws.Cells.Item(153, 4).Value = 0
$ws.Cells.Item(153, 5).Value = 11
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 0
$ws.Cells.Item(153, 8).Value = 0

# Row 154: 'Mali' -> 'Dominica'
$ws.Cells.Item(154, 1).Value = "Dominica"
$ws.Cells.Item(154, 2).Value = 11
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(154, 4).Value = 0
$ws.Cells.Item(154, 5).Value = 11
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 0

# Row 158: 'Surinam' -> 'Guinea'
$ws.Cells.Item(158, 1).Value = "Guinea"
$ws.Cells.Item(158, 2).Value = 8
$ws.Cells.Item(158, 3).Value = 4
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 5).Value = 8
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 0

# Row 159: 'Guinea' -> 'Haiti'
$ws.Cells.Item(159, 1).Value = "Haiti"
$ws.Cells.Item(159, 2).Value = 8
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 5).Value = 8
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 0

# Row 160: 'Haiti' -> 'Surinam'
$ws.Cells.Item(160, 1).Value = "Surinam"
$ws.Cells.Item(160, 2).Value = 8
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 0
$ws.Cells.Item(160, 5).Value = 8
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 0

# Row 169: 'Suazilandia' -> 'Laos'
$ws.Cells.Item(169, 1).Value = "Laos"
$ws.Cells.Item(169, 2).Value = 6
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(169, 4).Value = 0
$ws.Cells.Item(169, 5).Value = 6
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 0

# Row 170: 'Benin' -> 'Eritrea'
$ws.Cells.Item(170, 1).Value = "Eritrea"
$ws.Cells.Item(170, 2).Value = 6
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = 6
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 0

# Row 171: 'Laos' -> 'Benin'
$ws.Cells.Item(171, 1).Value = "Benin"
$ws.Cells.Item(171, 2).Value = 6
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 5).Value = 6
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 0

# Row 172: 'Eritrea' -> 'Suazilandia'
$ws.Cells.Item(172, 1).Value = "Suazilandia"
$ws.Cells.Item(172, 2).Value = 6
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 5).Value = 6
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 0

# Row 173: 'Siria' -> 'Fiyi'
$ws.Cells.Item(173, 1).Value = "Fiyi"
$ws.Cells.Item(173, 2).Value = 5
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 5
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

# Row 174: 'Birmania' -> 'Montserrat'
$ws.Cells.Item(174, 1).Value = "Montserrat"
$ws.Cells.Item(174, 2).Value = 5
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 5
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 0

# Row 175: 'Fiyi' -> 'Siria'
$ws.Cells.Item(175, 1).Value = "Siria"
$ws.Cells.Item(175, 2).Value = 5
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 0
$ws.Cells.Item(175, 5).Value = 5
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 0

# Row 176: 'Montserrat' -> 'Birmania'
$ws.Cells.Item(176, 1).Value = "Birmania"
$ws.Cells.Item(176, 2).Value = 5
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 0
$ws.Cells.Item(176, 5).Value = 5
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

# Row 177: 'Cabo Verde' -> 'Guyana'
$ws.Cells.Item(177, 1).Value = "Guyana"
$ws.Cells.Item(177, 2).Value = 5
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 5).Value = 4
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 1

# Row 178: 'Guyana' -> 'Cabo Verde'
$ws.Cells.Item(178, 1).Value = "Cabo Verde"
$ws.Cells.Item(178, 2).Value = 5
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 4).Value = 0
$ws.Cells.Item(178, 5).Value = 4
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(178, 8).Value = 1

# Row 180: 'Congo' -> 'Santa Sede'
$ws.Cells.Item(180, 1).Value = "Santa Sede"
$ws.Cells.Item(180, 2).Value = 4
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 4
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

# Row 181: 'Santa Sede' -> 'Congo'
$ws.Cells.Item(181, 1).Value = "Congo"
$ws.Cells.Item(181, 2).Value = 4
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 4
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 0

# Row 184: 'Mauritania' -> 'Liberia'
$ws.Cells.Item(184, 1).Value = "Liberia"
$ws.Cells.Item(184, 2).Value = 3
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 3
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

# Row 185: 'Republica de Africa Central' -> 'Mauritania'
$ws.Cells.Item(185, 1).Value = "Mauritania"
$ws.Cells.Item(185, 2).Value = 3
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 3
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

# Row 186: 'Liberia' -> 'Republica de Africa Central'
$ws.Cells.Item(186, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(186, 2).Value = 3
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 3
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

# Row 187: 'San Martin (Parte Holandesa)' -> 'Republica del Chad'
$ws.Cells.Item(187, 1).Value = "Republica del Chad"
$ws.Cells.Item(187, 2).Value = 3
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 3
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

# Row 188: 'Republica del Chad' -> 'San Martin (Parte Holandesa)'
$ws.Cells.Item(188, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(188, 2).Value = 3
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 5).Value = 3
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

# Row 192: 'Gambia' -> 'Santa Lucia'
$ws.Cells.Item(192, 1).Value = "Santa Lucia"
$ws.Cells.Item(192, 2).Value = 3
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 1
$ws.Cells.Item(192, 5).Value = 2
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 0

# Row 194: 'Santa Lucia' -> 'Gambia'
$ws.Cells.Item(194, 1).Value = "Gambia"
$ws.Cells.Item(194, 2).Value = 3
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 0
$ws.Cells.Item(194, 5).Value = 2
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 1

# Row 195: 'Islas Virgenes Britanicas' -> 'Belice'
$ws.Cells.Item(195, 1).Value = "Belice"
$ws.Cells.Item(195, 2).Value = 2
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 5).Value = 2
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

# Row 197: 'Guinea-Bisau' -> 'Anguila'
$ws.Cells.Item(197, 1).Value = "Anguila"
$ws.Cells.Item(197, 2).Value = 2
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 0
$ws.Cells.Item(197, 5).Value = 2
$ws.Cells.Item(197, 6).Value = 0
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 0

# Row 198: 'Anguila' -> 'Islas Turcas y Caicos'
$ws.Cells.Item(198, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(198, 2).Value = 2
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 0
$ws.Cells.Item(198, 5).Value = 2
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 0

# Row 199: 'Islas Turcas y Caicos' -> 'Guinea-Bisau'
$ws.Cells.Item(199, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(199, 2).Value = 2
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 0
$ws.Cells.Item(199, 5).Value = 2
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0

# Row 200: 'Belice' -> 'Islas Virgenes Britanicas'
$ws.Cells.Item(200, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(200, 2).Value = 2
$ws.Cells.Item(200, 3).Value = 0
$ws.Cells.Item(200, 4).Value = 0
$ws.Cells.Item(200, 5).Value = 2
$ws.Cells.Item(200, 6).Value = 0
$ws.Cells.Item(200, 7).Value = 0
$ws.Cells.Item(200, 8).Value = 0

# Row 202: 'San Vicente y las Granadinas' -> 'Papua Nueva Guinea'
$ws.Cells.Item(202, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(202, 2).Value = 1
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 0
$ws.Cells.Item(202, 5).Value = 1
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

# Row 203: 'Libia' -> 'Timor Oriental'
$ws.Cells.Item(203, 1).Value = "Timor Oriental"
$ws.Cells.Item(203, 2).Value = 1
$ws.Cells.Item(203, 3).Value = 0
$ws.Cells.Item(203, 4).Value = 0
$ws.Cells.Item(203, 5).Value = 1
$ws.Cells.Item(203, 6).Value = 0
$ws.Cells.Item(203, 7).Value = 0
$ws.Cells.Item(203, 8).Value = 0

# Row 204: 'Timor Oriental' -> 'Libia'
$ws.Cells.Item(204, 1).Value = "Libia"
$ws.Cells.Item(204, 2).Value = 1
$ws.Cells.Item(204, 3).Value = 0
$ws.Cells.Item(204, 4).Value = 0
$ws.Cells.Item(204, 5).Value = 1
$ws.Cells.Item(204, 6).Value = 0
$ws.Cells.Item(204, 7).Value = 0
$ws.Cells.Item(204, 8).Value = 0

# Row 205: 'Papua Nueva Guinea' -> 'San Vicente y las Granadinas'
$ws.Cells.Item(205, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(205, 2).Value = 1
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 0
$ws.Cells.Item(205, 5).Value = 1
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0
